# chore: update Sheets via scheduled runner
# Refresh market-price-derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the affected Leve rows across each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 866.5294
$ws.Range("I33").Value = 648.3077
$ws.Range("K33").Value = 648.3077
$ws.Range("M33").Value = -419.3077

$ws.Range("H92").Value = 666.1875
$ws.Range("I92").Value = 618.5
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 618.5
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 629.5
$ws.Range("N92").Value = -3496

$ws.Range("H121").Value = 1158.6
$ws.Range("J121").Value = 1191.7142
$ws.Range("L121").Value = 3575.1426
$ws.Range("N121").Value = -7069.142599999999

$ws.Range("H132").Value = 3391117
$ws.Range("I132").Value = 3922431.8
$ws.Range("K132").Value = 11767295.4
$ws.Range("M132").Value = -11764765.4

$ws.Range("H135").Value = 725.7727
$ws.Range("I135").Value = 793
$ws.Range("K135").Value = 7137
$ws.Range("M135").Value = -4602

$ws.Range("H137").Value = 3574721.2
$ws.Range("I137").Value = 5003369.5
$ws.Range("J137").Value = 3099.875
$ws.Range("K137").Value = 15010108.5
$ws.Range("L137").Value = 9299.625
$ws.Range("M137").Value = -15007558.5
$ws.Range("N137").Value = -14399.625

$ws.Range("H138").Value = 3618.1836
$ws.Range("I138").Value = 2064.2812
$ws.Range("J138").Value = 6543.1763
$ws.Range("K138").Value = 6192.8436
$ws.Range("L138").Value = 19629.5289
$ws.Range("M138").Value = -1052.8436
$ws.Range("N138").Value = -29909.5289

$ws.Range("H141").Value = 514237.62
$ws.Range("I141").Value = 1145
$ws.Range("J141").Value = 2669226.5
$ws.Range("K141").Value = 3435
$ws.Range("L141").Value = 8007679.5
$ws.Range("M141").Value = 1745
$ws.Range("N141").Value = -8018039.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3647.484
$ws.Range("I32").Value = 2367.8052
$ws.Range("K32").Value = 2367.8052
$ws.Range("M32").Value = -2080.8052

$ws.Range("H74").Value = 890.5454999999999
$ws.Range("I74").Value = 846.58826
$ws.Range("J74").Value = 1040
$ws.Range("K74").Value = 846.58826
$ws.Range("L74").Value = 1040
$ws.Range("M74").Value = 27.41174000000001
$ws.Range("N74").Value = -2788

$ws.Range("H77").Value = 890.5454999999999
$ws.Range("I77").Value = 846.58826
$ws.Range("J77").Value = 1040
$ws.Range("K77").Value = 4232.9413
$ws.Range("L77").Value = 5200
$ws.Range("M77").Value = 135.0586999999996
$ws.Range("N77").Value = -13936

$ws.Range("H97").Value = 858.125
$ws.Range("I97").Value = 838
$ws.Range("J97").Value = 999
$ws.Range("K97").Value = 838
$ws.Range("L97").Value = 999
$ws.Range("M97").Value = -342
$ws.Range("N97").Value = -1991

$ws.Range("H102").Value = 3739.8
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378

$ws.Range("H132").Value = 14495627
$ws.Range("I132").Value = 17244022
$ws.Range("J132").Value = 4090.7273
$ws.Range("K132").Value = 51732066
$ws.Range("L132").Value = 12272.1819
$ws.Range("M132").Value = -51729536
$ws.Range("N132").Value = -17332.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 70010
$ws.Range("J19").Value = 70010
$ws.Range("L19").Value = 70010
$ws.Range("N19").Value = -70356

$ws.Range("H49").Value = 30029.5
$ws.Range("I49").Value = 10059
$ws.Range("J49").Value = 50000
$ws.Range("K49").Value = 10059
$ws.Range("L49").Value = 50000
$ws.Range("M49").Value = -9820
$ws.Range("N49").Value = -50478

$ws.Range("H94").Value = 716.0909
$ws.Range("I94").Value = 508.55554
$ws.Range("J94").Value = 1650
$ws.Range("K94").Value = 508.55554
$ws.Range("L94").Value = 1650
$ws.Range("M94").Value = -57.55554000000001
$ws.Range("N94").Value = -2552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2162.7778
$ws.Range("I105").Value = 2074.4546
$ws.Range("J105").Value = 2301.5715
$ws.Range("K105").Value = 2074.4546
$ws.Range("L105").Value = 2301.5715
$ws.Range("M105").Value = -327.4546
$ws.Range("N105").Value = -5795.5715

$ws.Range("H107").Value = 1452.7273
$ws.Range("I107").Value = 301.8889
$ws.Range("J107").Value = 2833.7334
$ws.Range("K107").Value = 301.8889
$ws.Range("L107").Value = 2833.7334
$ws.Range("M107").Value = 1618.1111
$ws.Range("N107").Value = -6673.7334

$ws.Range("H132").Value = 3915.0715
$ws.Range("I132").Value = 2847.6365
$ws.Range("J132").Value = 4605.7646
$ws.Range("K132").Value = 8542.9095
$ws.Range("L132").Value = 13817.2938
$ws.Range("M132").Value = -6012.9095
$ws.Range("N132").Value = -18877.2938

$ws.Range("H134").Value = 1589.683
$ws.Range("I134").Value = 1230.8438
$ws.Range("J134").Value = 2865.5557
$ws.Range("K134").Value = 3692.5314
$ws.Range("L134").Value = 8596.667099999999
$ws.Range("M134").Value = -1157.5314
$ws.Range("N134").Value = -13666.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 1679.8667
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1679.8667
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 5039.6001
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -5319.6001

$ws.Range("H68").Value = 1868.5916
$ws.Range("J68").Value = 2490.5217
$ws.Range("L68").Value = 7471.5651
$ws.Range("N68").Value = -9093.5651

$ws.Range("H71").Value = 1868.5916
$ws.Range("J71").Value = 2490.5217
$ws.Range("L71").Value = 22414.6953
$ws.Range("N71").Value = -30526.6953

$ws.Range("H122").Value = 830.95
$ws.Range("I122").Value = 464.2857
$ws.Range("J122").Value = 1028.3846
$ws.Range("K122").Value = 4178.571300000001
$ws.Range("L122").Value = 9255.4614
$ws.Range("M122").Value = -1728.571300000001
$ws.Range("N122").Value = -14155.4614

$ws.Range("H131").Value = 1378
$ws.Range("J131").Value = 1231.1818
$ws.Range("L131").Value = 3693.5454
$ws.Range("N131").Value = -13773.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1659.3529
$ws.Range("I102").Value = 1098.909
$ws.Range("J102").Value = 2686.8333
$ws.Range("K102").Value = 1098.909
$ws.Range("L102").Value = 2686.8333
$ws.Range("M102").Value = 523.0909999999999
$ws.Range("N102").Value = -5930.8333

$ws.Range("H132").Value = 3385.3125
$ws.Range("I132").Value = 2494.9092
$ws.Range("J132").Value = 5344.2
$ws.Range("K132").Value = 7484.7276
$ws.Range("L132").Value = 16032.6
$ws.Range("M132").Value = -4954.7276
$ws.Range("N132").Value = -21092.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2379.3845
$ws.Range("I7").Value = 1585.3334
$ws.Range("J7").Value = 3060
$ws.Range("K7").Value = 1585.3334
$ws.Range("L7").Value = 3060
$ws.Range("M7").Value = -1473.3334
$ws.Range("N7").Value = -3284

$ws.Range("H20").Value = 19876.5
$ws.Range("J20").Value = 19876.5
$ws.Range("L20").Value = 19876.5
$ws.Range("N20").Value = -20328.5

$ws.Range("H40").Value = 2319.6
$ws.Range("I40").Value = 1697.25
$ws.Range("J40").Value = 2545.9092
$ws.Range("K40").Value = 1697.25
$ws.Range("L40").Value = 2545.9092
$ws.Range("M40").Value = -1561.25
$ws.Range("N40").Value = -2817.9092

$ws.Range("H46").Value = 1972.7273
$ws.Range("I46").Value = 920
$ws.Range("J46").Value = 2367.5
$ws.Range("K46").Value = 920
$ws.Range("L46").Value = 2367.5
$ws.Range("M46").Value = -732
$ws.Range("N46").Value = -2743.5

$ws.Range("H93").Value = 8500
$ws.Range("J93").Value = 8500
$ws.Range("L93").Value = 8500
$ws.Range("N93").Value = -10996

$ws.Range("H100").Value = 3200.7778
$ws.Range("I100").Value = 3333.3333
$ws.Range("J100").Value = 3134.5
$ws.Range("K100").Value = 3333.3333
$ws.Range("L100").Value = 3134.5
$ws.Range("M100").Value = -2792.3333
$ws.Range("N100").Value = -4216.5

$ws.Range("H122").Value = 2852.3572
$ws.Range("I122").Value = 2463.606
$ws.Range("J122").Value = 4277.778
$ws.Range("K122").Value = 7390.818000000001
$ws.Range("L122").Value = 12833.334
$ws.Range("M122").Value = -4940.818000000001
$ws.Range("N122").Value = -17733.334

$ws.Range("H126").Value = 2379.3845
$ws.Range("I126").Value = 1585.3334
$ws.Range("J126").Value = 3060
$ws.Range("K126").Value = 4756.0002
$ws.Range("L126").Value = 9180
$ws.Range("M126").Value = -2286.0002
$ws.Range("N126").Value = -14120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 39133.332
$ws.Range("J46").Value = 39133.332
$ws.Range("L46").Value = 39133.332
$ws.Range("N46").Value = -39595.332

$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 30000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -34340

$ws.Range("H134").Value = 39133.332
$ws.Range("J134").Value = 39133.332
$ws.Range("L134").Value = 117399.996
$ws.Range("N134").Value = -122469.996
